$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text to include "Rs." suffix
$ws.Range("G10").Value = "Vendor `nBasic Charge Rs."
$ws.Range("H10").Value = "Vendor `nTax Rs."

# Update the active selection to match the authored change
$ws.Range("G10").Select()
